# Insert a new "country" column between "address" (C) and "city" (D),
# shifting the existing city..is_published columns one to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Insert()

$ws.Range("D1").Value = "country"
$ws.Range("D2").Value = "India"

# Update the saved selection/view to match the post-edit state.
[void]$ws.Range("D3").Select()
